$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells that look like plain numbers keep their exact
# textual representation (e.g. trailing zeros) instead of being auto-converted
# to numeric values by Excel when assigned via .Value
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "59.900.02"
$ws.Range("E2").Value = "  +4.81%  "

$ws.Range("D3").Value = "3.336.19"
$ws.Range("E3").Value = "  +1.85%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "411.97"
$ws.Range("E5").Value = "  +3.08%  "

$ws.Range("D6").Value = "111.33"
$ws.Range("E6").Value = "  +0.86%  "

$ws.Range("D7").Value = "0.583"
$ws.Range("E7").Value = "  +4.33%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").Value = "0.630"
$ws.Range("E9").Value = "  +0.89%  "

$ws.Range("D10").Value = "39.73"
$ws.Range("E10").Value = "  +0.62%  "

$ws.Range("D11").Value = "0.0985"
$ws.Range("E11").Value = "  +3.14%  "

$ws.Range("D12").Value = "0.143"
$ws.Range("E12").Value = "  +1.14%  "

$ws.Range("D13").Value = "3.851.27"
$ws.Range("E13").Value = "  +1.95%  "

$ws.Range("D14").Value = "8.50"
$ws.Range("E14").Value = "  +4.83%  "

$ws.Range("D15").Value = "19.20"
$ws.Range("E15").Value = "  +0.27%  "

$ws.Range("D16").Value = "3.344.45"
$ws.Range("E16").Value = "  +2.17%  "

$ws.Range("E17").Value = "  -0.44%  "

$ws.Range("D18").Value = "59.604.19"
$ws.Range("E18").Value = "  +4.76%  "

$ws.Range("D19").Value = "10.63"
$ws.Range("E19").Value = "  -2.60%  "

$ws.Range("E20").Value = "  +0.54%  "

$ws.Range("E21").Value = "  +3.09%  "

$ws.Range("D22").Value = "13.01"
$ws.Range("E22").Value = "  +0.67%  "

$ws.Range("D23").Value = "302.06"
$ws.Range("E23").Value = "  -0.89%  "

$ws.Range("D24").Value = "75.25"
$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("D25").Value = "3.19"
$ws.Range("E25").Value = "  -0.16%  "

$ws.Range("D26").Value = "0.184"
$ws.Range("E26").Value = "  +8.43%  "

$ws.Range("D27").Value = "28.53"
$ws.Range("E27").Value = "  +0.57%  "

$ws.Range("E28").Value = "  +1.87%  "

$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "7.83"
$ws.Range("E29").Value = "  -2.39%  "

$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "7.43"
$ws.Range("E30").Value = "  +1.85%  "

$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "2.67"
$ws.Range("E31").Value = "  +24.97%  "

$ws.Range("D32").Value = "0.115"
$ws.Range("E32").Value = "  +3.44%  "

$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("D34").Value = "11.53"
$ws.Range("E34").Value = "  +4.31%  "

$ws.Range("D35").Value = "39.51"
$ws.Range("E35").Value = "  +5.15%  "

$ws.Range("D36").Value = "0.0503"
$ws.Range("E36").Value = "  +3.51%  "

$ws.Range("D37").Value = "51.70"
$ws.Range("E37").Value = "  +0.18%  "

$ws.Range("D38").Value = "3.16"
$ws.Range("E38").Value = "  -0.77%  "

$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  +0.14%  "

$ws.Range("E40").Value = "  -5.14%  "

$ws.Range("D41").Value = "138.66"
$ws.Range("E41").Value = "  +3.52%  "

$ws.Range("E42").Value = "  +2.28%  "

$ws.Range("E43").Value = "  -1.22%  "

$ws.Range("D44").Value = "0.284"
$ws.Range("E44").Value = "  +0.84%  "

$ws.Range("E45").Value = "  -2.30%  "

$ws.Range("D46").Value = "16.78"
$ws.Range("E46").Value = "  -4.43%  "

$ws.Range("E47").Value = "  +8.70%  "

$ws.Range("D48").Value = "22.29"
$ws.Range("E48").Value = "  +0.58%  "

$ws.Range("D49").Value = "2.195.01"
$ws.Range("E49").Value = "  +1.93%  "

$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Value = "2.13"
$ws.Range("E50").Value = "  +4.97%  "

$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").Value = "2.39"
$ws.Range("E51").Value = "  -0.05%  "
